$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction4")

# Clear the old row contents (A1:N1) so the used range shrinks down
$ws.Range("A1:N1").ClearContents()

# Write the new, smaller dataset
$ws.Range("A1").Value = 8
$ws.Range("B1").Value = 9
